$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 777.63
$ws.Range("I15").Value = 777.63
$ws.Range("K15").Value = 2332.89
$ws.Range("M15").Value = -2163.89

$ws.Range("H28").Value = 935.9048
$ws.Range("I28").Value = 926.3333
$ws.Range("J28").Value = 993.3333
$ws.Range("K28").Value = 926.3333
$ws.Range("L28").Value = 993.3333
$ws.Range("M28").Value = -441.3333
$ws.Range("N28").Value = -1963.3333

$ws.Range("H41").Value = 308.53845
$ws.Range("I41").Value = 254.55556
$ws.Range("J41").Value = 430
$ws.Range("K41").Value = 254.55556
$ws.Range("L41").Value = 430
$ws.Range("M41").Value = 185.44444
$ws.Range("N41").Value = -1310

$ws.Range("H86").Value = 3220.7778
$ws.Range("I86").Value = 1916.409
$ws.Range("K86").Value = 1916.409
$ws.Range("M86").Value = -793.4090000000001

$ws.Range("H89").Value = 3220.7778
$ws.Range("I89").Value = 1916.409
$ws.Range("K89").Value = 9582.045
$ws.Range("M89").Value = -3966.045

$ws.Range("H98").Value = 1856.8541
$ws.Range("I98").Value = 2098.4878
$ws.Range("J98").Value = 441.57144
$ws.Range("K98").Value = 2098.4878
$ws.Range("L98").Value = 441.57144
$ws.Range("M98").Value = -600.4877999999999
$ws.Range("N98").Value = -3437.57144

$ws.Range("H111").Value = 610
$ws.Range("I111").Value = 596.5
$ws.Range("J111").Value = 650.5
$ws.Range("K111").Value = 1789.5
$ws.Range("L111").Value = 1951.5
$ws.Range("M111").Value = 1277.5
$ws.Range("N111").Value = -8085.5

$ws.Range("H116").Value = 2854.7
$ws.Range("I116").Value = 2297.1052
$ws.Range("J116").Value = 3817.818
$ws.Range("K116").Value = 2297.1052
$ws.Range("L116").Value = 3817.818
$ws.Range("M116").Value = 1144.8948
$ws.Range("N116").Value = -10701.818

$ws.Range("H122").Value = 1856.8541
$ws.Range("I122").Value = 2098.4878
$ws.Range("J122").Value = 441.57144
$ws.Range("K122").Value = 6295.4634
$ws.Range("L122").Value = 1324.71432
$ws.Range("M122").Value = -3845.4634
$ws.Range("N122").Value = -6224.71432

$ws.Range("H131").Value = 1987.7778
$ws.Range("I131").Value = 1783.5714
$ws.Range("J131").Value = 2702.5
$ws.Range("K131").Value = 5350.7142
$ws.Range("L131").Value = 8107.5
$ws.Range("M131").Value = -310.7142000000003
$ws.Range("N131").Value = -18187.5

$ws.Range("H132").Value = 315818.97
$ws.Range("I132").Value = 325845.38
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 977536.14
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -975006.14
$ws.Range("N132").Value = -20060

$ws.Range("H137").Value = 27779662
$ws.Range("I137").Value = 1286.6428
$ws.Range("K137").Value = 3859.9284
$ws.Range("M137").Value = -1309.9284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11834.47
$ws.Range("I32").Value = 8788.902
$ws.Range("K32").Value = 8788.902
$ws.Range("M32").Value = -8501.902

$ws.Range("H45").Value = 907.7
$ws.Range("I45").Value = 884.625
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 884.625
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -507.625
$ws.Range("N45").Value = -1754

$ws.Range("H122").Value = 1163.6818
$ws.Range("I122").Value = 1173.1
$ws.Range("J122").Value = 1069.5
$ws.Range("K122").Value = 3519.3
$ws.Range("L122").Value = 3208.5
$ws.Range("M122").Value = -1069.3
$ws.Range("N122").Value = -8108.5

$ws.Range("H132").Value = 3462.9473
$ws.Range("I132").Value = 2669.3333
$ws.Range("J132").Value = 4177.2
$ws.Range("K132").Value = 8007.999899999999
$ws.Range("L132").Value = 12531.6
$ws.Range("M132").Value = -5477.999899999999
$ws.Range("N132").Value = -17591.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2250
$ws.Range("I107").Value = 2225
$ws.Range("J107").Value = 2275
$ws.Range("K107").Value = 2225
$ws.Range("L107").Value = 2275
$ws.Range("M107").Value = -305
$ws.Range("N107").Value = -6115

$ws.Range("H135").Value = 28674.4
$ws.Range("J135").Value = 28674.4
$ws.Range("L135").Value = 28674.4
$ws.Range("N135").Value = -38814.4

$ws.Range("H141").Value = 50513
$ws.Range("I141").Value = 75000
$ws.Range("J141").Value = 45615.6
$ws.Range("K141").Value = 75000
$ws.Range("L141").Value = 45615.6
$ws.Range("M141").Value = -69820
$ws.Range("N141").Value = -55975.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2835.2942
$ws.Range("I16").Value = 3014.2856
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 3014.2856
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -2727.2856
$ws.Range("N16").Value = -2574

$ws.Range("H31").Value = 1389.6086
$ws.Range("I31").Value = 1313.7222
$ws.Range("J31").Value = 1662.8
$ws.Range("K31").Value = 1313.7222
$ws.Range("L31").Value = 1662.8
$ws.Range("M31").Value = -1018.7222
$ws.Range("N31").Value = -2252.8

$ws.Range("H34").Value = 1389.6086
$ws.Range("I34").Value = 1313.7222
$ws.Range("J34").Value = 1662.8
$ws.Range("K34").Value = 1313.7222
$ws.Range("L34").Value = 1662.8
$ws.Range("M34").Value = -1111.7222
$ws.Range("N34").Value = -2066.8

$ws.Range("H105").Value = 697.1429000000001
$ws.Range("I105").Value = 575
$ws.Range("J105").Value = 860
$ws.Range("K105").Value = 575
$ws.Range("L105").Value = 860
$ws.Range("M105").Value = 1172
$ws.Range("N105").Value = -4354

$ws.Range("H113").Value = 2835.2942
$ws.Range("I113").Value = 3014.2856
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 3014.2856
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -844.2856000000002
$ws.Range("N113").Value = -6340

$ws.Range("H122").Value = 1163.1428
$ws.Range("I122").Value = 994
$ws.Range("J122").Value = 1520.2222
$ws.Range("K122").Value = 2982
$ws.Range("L122").Value = 4560.6666
$ws.Range("M122").Value = -532
$ws.Range("N122").Value = -9460.6666

$ws.Range("H132").Value = 2596.842
$ws.Range("I132").Value = 1591.1111
$ws.Range("J132").Value = 3502
$ws.Range("K132").Value = 4773.3333
$ws.Range("L132").Value = 10506
$ws.Range("M132").Value = -2243.3333
$ws.Range("N132").Value = -15566

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 646.9167
$ws.Range("I113").Value = 612.1818
$ws.Range("J113").Value = 676.3077
$ws.Range("K113").Value = 1836.5454
$ws.Range("L113").Value = 2028.9231
$ws.Range("M113").Value = 333.4546
$ws.Range("N113").Value = -6368.9231

$ws.Range("H140").Value = 1047.25
$ws.Range("I140").Value = 1047.25
$ws.Range("K140").Value = 3141.75
$ws.Range("M140").Value = 2038.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7782.1875
$ws.Range("I122").Value = 9125.875
$ws.Range("J122").Value = 6438.5
$ws.Range("K122").Value = 27377.625
$ws.Range("L122").Value = 19315.5
$ws.Range("M122").Value = -24927.625
$ws.Range("N122").Value = -24215.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 304.63635
$ws.Range("I22").Value = 262.625
$ws.Range("J22").Value = 416.66666
$ws.Range("K22").Value = 262.625
$ws.Range("L22").Value = 416.66666
$ws.Range("M22").Value = 32.375
$ws.Range("N22").Value = -1006.66666

$ws.Range("H27").Value = 304.63635
$ws.Range("I27").Value = 262.625
$ws.Range("J27").Value = 416.66666
$ws.Range("K27").Value = 262.625
$ws.Range("L27").Value = 416.66666
$ws.Range("M27").Value = -155.625
$ws.Range("N27").Value = -630.66666

$ws.Range("H40").Value = 2461.6667
$ws.Range("I40").Value = 2499.25
$ws.Range("J40").Value = 2418.7144
$ws.Range("K40").Value = 2499.25
$ws.Range("L40").Value = 2418.7144
$ws.Range("M40").Value = -2363.25
$ws.Range("N40").Value = -2690.7144

$ws.Range("H82").Value = 1816.9166
$ws.Range("I82").Value = 1475
$ws.Range("J82").Value = 2500.75
$ws.Range("K82").Value = 1475
$ws.Range("L82").Value = 2500.75
$ws.Range("M82").Value = -1114
$ws.Range("N82").Value = -3222.75

$ws.Range("H85").Value = 1816.9166
$ws.Range("I85").Value = 1475
$ws.Range("J85").Value = 2500.75
$ws.Range("K85").Value = 1475
$ws.Range("L85").Value = 2500.75
$ws.Range("M85").Value = -227
$ws.Range("N85").Value = -4996.75

$ws.Range("H122").Value = 17439
$ws.Range("I122").Value = 23402.4
$ws.Range("K122").Value = 70207.20000000001
$ws.Range("M122").Value = -67757.20000000001

$ws.Range("H132").Value = 6754.1665
$ws.Range("I132").Value = 7829.5293
$ws.Range("K132").Value = 23488.5879
$ws.Range("M132").Value = -20958.5879

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8714.777
$ws.Range("I136").Value = 10790.667
$ws.Range("K136").Value = 32372.001
$ws.Range("M136").Value = -29822.001
